{"js": "// Fix the typo \"Perquisite\" -> \"Prerequisite\" in the second paragraph\n// (the sentence describing actors/pre-condition/post-condition for the\n// \"make sale\" use case). This is the single textual change behind the\n// commit \"Sale Order Line done\" \u2014 the rest of the diff is just Word\n// re-splitting runs around the edited word, which does not change the\n// rendered text.\nconst body = context.document.body;\n\nconst results = body.search(\"Perquisite\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (const result of results.items) {\n  result.insertText(\"Prerequisite\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Fix the typo \"Perquisite\" -> \"Prerequisite\" in the paragraph describing\n# the actors/pre-condition/post-condition for the \"make sale\" use case.\n# This is the single textual change behind the commit \"Sale Order Line done\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$find.Execute(\n    \"Perquisite\",\n    $true,\n    $true,\n    $false,\n    $false,\n    $false,\n    $true,\n    $wdFindContinue,\n    $false,\n    \"Prerequisite\",\n    $wdReplaceOne\n)\n"}
